$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new row of data, copying row 10's formatting as a template
# (same column styling pattern: B bold/wrap, D vertical-top, E yellow fill).
$ws.Range("B10:E10").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)

$ws.Range("B11").Value = "GFG"
$ws.Range("C11").Value = "Top View of Binary Tree"
$ws.Range("D11").Value = "Java/Python"
$ws.Range("E11").Value = "Medium"

# Move the active selection to C17 (shifted down by the new row), matching
# the sheetView selection recorded in the saved workbook.
$ws.Range("C17").Select()
